$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1728.9474
$ws.Range("J40").Value = 2662.5
$ws.Range("L40").Value = 2662.5
$ws.Range("N40").Value = -3012.5

$ws.Range("H62").Value = 3113.5293
$ws.Range("I62").Value = 2337.9
$ws.Range("K62").Value = 2337.9
$ws.Range("M62").Value = -1713.9

$ws.Range("H65").Value = 3113.5293
$ws.Range("I65").Value = 2337.9
$ws.Range("K65").Value = 11689.5
$ws.Range("M65").Value = -8569.5

$ws.Range("H86").Value = 8067.579
$ws.Range("I86").Value = 12410
$ws.Range("J86").Value = 3242.6667
$ws.Range("K86").Value = 12410
$ws.Range("L86").Value = 3242.6667
$ws.Range("M86").Value = -11287
$ws.Range("N86").Value = -5488.6667

$ws.Range("H89").Value = 8067.579
$ws.Range("I89").Value = 12410
$ws.Range("J89").Value = 3242.6667
$ws.Range("K89").Value = 62050
$ws.Range("L89").Value = 16213.3335
$ws.Range("M89").Value = -56434
$ws.Range("N89").Value = -27445.3335

$ws.Range("H116").Value = 3369.7026
$ws.Range("I116").Value = 2587.6538
$ws.Range("J116").Value = 5218.1816
$ws.Range("K116").Value = 2587.6538
$ws.Range("L116").Value = 5218.1816
$ws.Range("M116").Value = 854.3462
$ws.Range("N116").Value = -12102.1816

$ws.Range("H121").Value = 1368.7142
$ws.Range("I121").Value = 520.75
$ws.Range("J121").Value = 2499.3333
$ws.Range("K121").Value = 1562.25
$ws.Range("L121").Value = 7497.999899999999
$ws.Range("M121").Value = 184.75
$ws.Range("N121").Value = -10991.9999

$ws.Range("H132").Value = 7149590
$ws.Range("I132").Value = 9097909
$ws.Range("K132").Value = 27293727
$ws.Range("M132").Value = -27291197

$ws.Range("H141").Value = 886186.75
$ws.Range("I141").Value = 1858.909
$ws.Range("J141").Value = 3665502.8
$ws.Range("K141").Value = 5576.727000000001
$ws.Range("L141").Value = 10996508.4
$ws.Range("M141").Value = -396.7270000000008
$ws.Range("N141").Value = -11006868.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1680.1
$ws.Range("I110").Value = 662.625
$ws.Range("J110").Value = 5750
$ws.Range("K110").Value = 662.625
$ws.Range("L110").Value = 5750
$ws.Range("M110").Value = 1382.375
$ws.Range("N110").Value = -9840

$ws.Range("H132").Value = 25004498
$ws.Range("I132").Value = 38465908
$ws.Range("J132").Value = 4737.5713
$ws.Range("K132").Value = 115397724
$ws.Range("L132").Value = 14212.7139
$ws.Range("M132").Value = -115395194
$ws.Range("N132").Value = -19272.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1996.875
$ws.Range("I86").Value = 1380.7693
$ws.Range("J86").Value = 4666.6665
$ws.Range("K86").Value = 1380.7693
$ws.Range("L86").Value = 4666.6665
$ws.Range("M86").Value = -257.7692999999999
$ws.Range("N86").Value = -6912.6665

$ws.Range("H89").Value = 1996.875
$ws.Range("I89").Value = 1380.7693
$ws.Range("J89").Value = 4666.6665
$ws.Range("K89").Value = 6903.8465
$ws.Range("L89").Value = 23333.3325
$ws.Range("M89").Value = -1287.8465
$ws.Range("N89").Value = -34565.3325

$ws.Range("H99").Value = 1966
$ws.Range("I99").Value = 1232.2222
$ws.Range("J99").Value = 3066.6667
$ws.Range("K99").Value = 1232.2222
$ws.Range("L99").Value = 3066.6667
$ws.Range("M99").Value = 265.7778000000001
$ws.Range("N99").Value = -6062.6667

$ws.Range("H134").Value = 3745.1924
$ws.Range("I134").Value = 3653.75
$ws.Range("K134").Value = 10961.25
$ws.Range("M134").Value = -8426.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 70011
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 70011
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 70011
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -70351

$ws.Range("H41").Value = 3976.3333
$ws.Range("I41").Value = 1473.25
$ws.Range("J41").Value = 8982.5
$ws.Range("K41").Value = 1473.25
$ws.Range("L41").Value = 8982.5
$ws.Range("M41").Value = -1045.25
$ws.Range("N41").Value = -9838.5

$ws.Range("H50").Value = 14900
$ws.Range("J50").Value = 14900
$ws.Range("L50").Value = 14900
$ws.Range("N50").Value = -16150

$ws.Range("H51").Value = 166679760
$ws.Range("I51").Value = 1000000000
$ws.Range("J51").Value = 15720
$ws.Range("K51").Value = 1000000000
$ws.Range("L51").Value = 15720
$ws.Range("M51").Value = -999999264
$ws.Range("N51").Value = -17192

$ws.Range("H60").Value = 11857.857
$ws.Range("I60").Value = 14000
$ws.Range("J60").Value = 11693.077
$ws.Range("K60").Value = 14000
$ws.Range("L60").Value = 11693.077
$ws.Range("M60").Value = -13489
$ws.Range("N60").Value = -12715.077

$ws.Range("H61").Value = 166679760
$ws.Range("I61").Value = 1000000000
$ws.Range("J61").Value = 15720
$ws.Range("K61").Value = 1000000000
$ws.Range("L61").Value = 15720
$ws.Range("M61").Value = -999999652
$ws.Range("N61").Value = -16416

$ws.Range("H132").Value = 4477.913
$ws.Range("I132").Value = 3308.6155
$ws.Range("K132").Value = 9925.8465
$ws.Range("M132").Value = -7395.8465

$ws.Range("H134").Value = 1552.5714
$ws.Range("I134").Value = 949.5789
$ws.Range("J134").Value = 2825.5557
$ws.Range("K134").Value = 2848.7367
$ws.Range("L134").Value = 8476.667099999999
$ws.Range("M134").Value = -313.7366999999999
$ws.Range("N134").Value = -13546.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 87.875
$ws.Range("I2").Value = 34.5
$ws.Range("J2").Value = 105.666664
$ws.Range("K2").Value = 207
$ws.Range("L2").Value = 633.999984
$ws.Range("M2").Value = -94
$ws.Range("N2").Value = -859.999984

$ws.Range("H5").Value = 1641.5625
$ws.Range("I5").Value = 547.0833
$ws.Range("J5").Value = 4925
$ws.Range("K5").Value = 1641.2499
$ws.Range("L5").Value = 14775
$ws.Range("M5").Value = -1529.2499
$ws.Range("N5").Value = -14999

$ws.Range("H38").Value = 126.42857
$ws.Range("I38").Value = 97.5
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 292.5
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = 54.5
$ws.Range("N38").Value = -1594

$ws.Range("H131").Value = 43066.812
$ws.Range("I131").Value = 2311.4285
$ws.Range("J131").Value = 50991.473
$ws.Range("K131").Value = 6934.2855
$ws.Range("L131").Value = 152974.419
$ws.Range("M131").Value = -1894.2855
$ws.Range("N131").Value = -163054.419

$ws.Range("H134").Value = 2202.2354
$ws.Range("I134").Value = 970.9
$ws.Range("J134").Value = 3961.2856
$ws.Range("K134").Value = 2912.7
$ws.Range("L134").Value = 11883.8568
$ws.Range("M134").Value = 2157.3
$ws.Range("N134").Value = -22023.8568

$ws.Range("H135").Value = 1641.5625
$ws.Range("I135").Value = 547.0833
$ws.Range("J135").Value = 4925
$ws.Range("K135").Value = 4923.7497
$ws.Range("L135").Value = 44325
$ws.Range("M135").Value = -2388.7497
$ws.Range("N135").Value = -49395

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 46113.2
$ws.Range("J4").Value = 46113.2
$ws.Range("L4").Value = 46113.2
$ws.Range("N4").Value = -46337.2

$ws.Range("H26").Value = 30014
$ws.Range("J26").Value = 30014
$ws.Range("L26").Value = 30014
$ws.Range("N26").Value = -30574

$ws.Range("H50").Value = 30014
$ws.Range("J50").Value = 30014
$ws.Range("L50").Value = 30014
$ws.Range("N50").Value = -31010

$ws.Range("H132").Value = 3164.027
$ws.Range("I132").Value = 2713.55
$ws.Range("J132").Value = 3694
$ws.Range("K132").Value = 8140.650000000001
$ws.Range("L132").Value = 11082
$ws.Range("M132").Value = -5610.650000000001
$ws.Range("N132").Value = -16142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 71430510
$ws.Range("I61").Value = 76924390
$ws.Range("K61").Value = 76924390
$ws.Range("M61").Value = -76924188

$ws.Range("H113").Value = 71430510
$ws.Range("I113").Value = 76924390
$ws.Range("K113").Value = 76924390
$ws.Range("M113").Value = -76922220

$ws.Range("H132").Value = 3130.1
$ws.Range("I132").Value = 1892.1666
$ws.Range("J132").Value = 3955.389
$ws.Range("K132").Value = 5676.4998
$ws.Range("L132").Value = 11866.167
$ws.Range("M132").Value = -3146.4998
$ws.Range("N132").Value = -16926.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2999.5
$ws.Range("I62").Value = 2999
$ws.Range("K62").Value = 2999
$ws.Range("M62").Value = -2375

$ws.Range("H65").Value = 2999.5
$ws.Range("I65").Value = 2999
$ws.Range("K65").Value = 14995
$ws.Range("M65").Value = -11875

$ws.Range("H123").Value = 21134.5
$ws.Range("J123").Value = 21134.5
$ws.Range("L123").Value = 21134.5
$ws.Range("N123").Value = -30934.5

$ws.Range("H132").Value = 226937.31
$ws.Range("I132").Value = 279824.47
$ws.Range("J132").Value = 15388.667
$ws.Range("K132").Value = 839473.4099999999
$ws.Range("L132").Value = 46166.001
$ws.Range("M132").Value = -836943.4099999999
$ws.Range("N132").Value = -51226.001
